$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-4 (columns D-O) ---
$ws.Range("D2").Value = 0.0774398143281471
$ws.Range("E2").Value = 0.0300215734496843
$ws.Range("F2").Value = 1.08051719884797
$ws.Range("G2").Value = 1.01877154148646
$ws.Range("H2").Value = 1.14600513408804
$ws.Range("I2").Value = 2.57947220714247
$ws.Range("J2").Value = 0.0100178771177105
$ws.Range("K2").Value = 0.181666500331872
$ws.Range("L2").Value = 0.1717473063965
$ws.Range("M2").Value = 2423
$ws.Range("N2").Value = 1170
$ws.Range("O2").Value = 51.7127527858027
$ws.Range("D3").Value = -0.0837455850970771
$ws.Range("E3").Value = 0.0301465959358952
$ws.Range("F3").Value = 0.919665202867399
$ws.Range("G3").Value = 0.866898904218383
$ws.Range("H3").Value = 0.975643274261275
$ws.Range("I3").Value = -2.77794498838796
$ws.Range("J3").Value = 0.00555927287493991
$ws.Range("K3").Value = 0.1816830876367
$ws.Range("L3").Value = 0.171729623889449
$ws.Range("M3").Value = 2423
$ws.Range("N3").Value = 1166
$ws.Range("O3").Value = 51.8778373916632
$ws.Range("D4").Value = -0.0437922762272516
$ws.Range("E4").Value = 0.0291651820475902
$ws.Range("F4").Value = 0.957152760206752
$ws.Range("G4").Value = 0.903972776946592
$ws.Range("H4").Value = 1.01346127863044
$ws.Range("I4").Value = -1.50152590015703
$ws.Range("J4").Value = 0.133512163509766
$ws.Range("K4").Value = 0.179350738724407
$ws.Range("L4").Value = 0.168702833573473
$ws.Range("M4").Value = 2423
$ws.Range("N4").Value = 1094
$ws.Range("O4").Value = 54.8493602971523

# --- Add new rows 5-10 (columns A-O) ---
$ws.Range("A5").Value = "AEDB.CEA"
$ws.Range("B5").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C5").Value = "Macrophages_rank"
$ws.Range("D5").Value = 0.0868178343037033
$ws.Range("E5").Value = 0.0285581748232536
$ws.Range("F5").Value = 1.09069797386427
$ws.Range("G5").Value = 1.03132441044854
$ws.Range("H5").Value = 1.1534896858247
$ws.Range("I5").Value = 3.04003441540009
$ws.Range("J5").Value = 0.00241871782439107
$ws.Range("K5").Value = 0.247483934421806
$ws.Range("L5").Value = 0.238370418056672
$ws.Range("M5").Value = 2423
$ws.Range("N5").Value = 1171
$ws.Range("O5").Value = 51.6714816343376
$ws.Range("A6").Value = "AEDB.CEA"
$ws.Range("B6").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C6").Value = "SMC_rank"
$ws.Range("D6").Value = -0.148492575381802
$ws.Range("E6").Value = 0.0284161560212707
$ws.Range("F6").Value = 0.862006407215198
$ws.Range("G6").Value = 0.81530888203239
$ws.Range("H6").Value = 0.911378573759403
$ws.Range("I6").Value = -5.22563907907348
$ws.Range("J6").Value = 0.000000205891830271717
$ws.Range("K6").Value = 0.258832789468056
$ws.Range("L6").Value = 0.249825549062286
$ws.Range("M6").Value = 2423
$ws.Range("N6").Value = 1167
$ws.Range("O6").Value = 51.8365662401981
$ws.Range("A7").Value = "AEDB.CEA"
$ws.Range("B7").Value = "MCP1_pg_ml_2015_rank"
$ws.Range("C7").Value = "VesselDensity_rank"
$ws.Range("D7").Value = -0.0335025123311633
$ws.Range("E7").Value = 0.0278112126946543
$ws.Range("F7").Value = 0.967052481672265
$ws.Range("G7").Value = 0.915749439662387
$ws.Range("H7").Value = 1.02122967462941
$ws.Range("I7").Value = -1.2046404699786
$ws.Range("J7").Value = 0.228606040036961
$ws.Range("K7").Value = 0.237564890802019
$ws.Range("L7").Value = 0.227681472719822
$ws.Range("M7").Value = 2423
$ws.Range("N7").Value = 1095
$ws.Range("O7").Value = 54.8080891456872
$ws.Range("A8").Value = "AEDB.CEA"
$ws.Range("B8").Value = "MCP1_rank"
$ws.Range("C8").Value = "Macrophages_rank"
$ws.Range("D8").Value = 0.116278510133553
$ws.Range("E8").Value = 0.0387525820093411
$ws.Range("F8").Value = 1.12330868142184
$ws.Range("G8").Value = 1.04114747825829
$ws.Range("H8").Value = 1.21195356095809
$ws.Range("I8").Value = 3.00053581218212
$ws.Range("J8").Value = 0.0028179014628275
$ws.Range("K8").Value = 0.0898436194806476
$ws.Range("L8").Value = 0.0781962800590106
$ws.Range("M8").Value = 2423
$ws.Range("N8").Value = 555
$ws.Range("O8").Value = 77.0945109368551
$ws.Range("A9").Value = "AEDB.CEA"
$ws.Range("B9").Value = "MCP1_rank"
$ws.Range("C9").Value = "SMC_rank"
$ws.Range("D9").Value = -0.251422760451051
$ws.Range("E9").Value = 0.0404251520664834
$ws.Range("F9").Value = 0.777693523987001
$ws.Range("G9").Value = 0.718452232087924
$ws.Range("H9").Value = 0.841819664883863
$ws.Range("I9").Value = -6.21946356658252
$ws.Range("J9").Value = 0.000000000995217476255572
$ws.Range("K9").Value = 0.138234731950483
$ws.Range("L9").Value = 0.12714584063367
$ws.Range("M9").Value = 2423
$ws.Range("N9").Value = 552
$ws.Range("O9").Value = 77.2183243912505
$ws.Range("A10").Value = "AEDB.CEA"
$ws.Range("B10").Value = "MCP1_rank"
$ws.Range("C10").Value = "VesselDensity_rank"
$ws.Range("D10").Value = -0.0465296995741358
$ws.Range("E10").Value = 0.0510257920998464
$ws.Range("F10").Value = 0.954536210830026
$ws.Range("G10").Value = 0.8636909662687
$ws.Range("H10").Value = 1.0549367926378
$ws.Range("I10").Value = -0.911885884751915
$ws.Range("J10").Value = 0.362238728725745
$ws.Range("K10").Value = 0.0830549675526182
$ws.Range("L10").Value = 0.0710799391437904
$ws.Range("M10").Value = 2423
$ws.Range("N10").Value = 544
$ws.Range("O10").Value = 77.5484936029715
